$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the 4 new rows needed (from bottom to top so row indices stay valid)
$ws.Rows.Item(19).Insert()   # room for document_chunks.updated_at (before status)
$ws.Rows.Item(12).Insert()   # room for documents.updated_at (after created_at)
$ws.Rows.Item(6).Insert()    # room for users.updated_at (after created_at)
$ws.Rows.Item(5).Insert()    # room for users.is_admin (before created_at)

# Fill in the new rows' data: users.is_admin (row 5)
$ws.Range("A5").Value = "users"
$ws.Range("B5").Value = "is_admin"
$ws.Range("C5").Value = "Boolean"
$ws.Range("D5").Value = "Default False"
$ws.Range("E5").Value = "Is admin user"

# users.updated_at (row 7)
$ws.Range("A7").Value = "users"
$ws.Range("B7").Value = "updated_at"
$ws.Range("C7").Value = "DateTime"
$ws.Range("D7").Value = "Not Null"
$ws.Range("E7").Value = "User update timestamp"

# documents.updated_at (row 14)
$ws.Range("A14").Value = "documents"
$ws.Range("B14").Value = "updated_at"
$ws.Range("C14").Value = "DateTime"
$ws.Range("D14").Value = "Not Null"
$ws.Range("E14").Value = "Document update timestamp"

# document_chunks.updated_at (row 22)
$ws.Range("A22").Value = "document_chunks"
$ws.Range("B22").Value = "updated_at"
$ws.Range("C22").Value = "DateTime"
$ws.Range("D22").Value = "Not Null"
$ws.Range("E22").Value = "Chunk update timestamp"

# Update selection to match the target workbook
$ws.Range("E2").Select() | Out-Null
